$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: 'Bitcoin'
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '45.661.39'
$ws.Range('E2').Value = '  +6.08%  '

# Row 3: 'Ethereum'
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.399.46'
$ws.Range('E3').Value = '  +4.15%  '

# Row 4: 'TetherUSD'
$ws.Range('E4').Value = '  -0.14%  '

# Row 5: 'Solana'
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '115.43'
$ws.Range('E5').Value = '  +10.24%  '

# Row 6: 'BNB'
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '320.12'
$ws.Range('E6').Value = '  +3.09%  '

# Row 7: 'XRP'
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.638'
$ws.Range('E7').Value = '  +2.04%  '

# Row 8: 'USDC'
$ws.Range('E8').Value = '  -0.41%  '

# Row 9: 'Cardano'
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.630'
$ws.Range('E9').Value = '  +3.96%  '

# Row 10: 'Avalanche'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '42.79'
$ws.Range('E10').Value = '  +7.91%  '

# Row 11: 'Dogecoin'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0934'
$ws.Range('E11').Value = '  +3.39%  '

# Row 12: 'Polkadot'
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '8.70'
$ws.Range('E12').Value = '  +4.98%  '

# Row 13: 'TRON'
$ws.Range('E13').Value = '  +3.06%  '

# Row 14: 'Polygon'
$ws.Range('E14').Value = '  +2.69%  '

# Row 15: 'Chainlink'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.00'
$ws.Range('E15').Value = '  +4.10%  '

# Row 16: 'WrappedliquidstakedEther2.0'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.764.49'
$ws.Range('E16').Value = '  -0.69%  '

# Row 17: 'WrappedEther'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.404.13'
$ws.Range('E17').Value = '  +4.60%  '

# Row 18: 'WrappedBTC'
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '45.674.85'
$ws.Range('E18').Value = '  +6.50%  '

# Row 19: 'Uniswap'
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '7.53'
$ws.Range('E19').Value = '  +2.78%  '

# Row 20: 'ShibaInu'
$ws.Range('E20').Value = '  +3.60%  '

# Row 21: 'InternetComputer(DFINITY)'
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '13.57'
$ws.Range('E21').Value = '  +0.24%  '

# Row 22: 'Litecoin'
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '75.18'
$ws.Range('E22').Value = '  +2.44%  '

# Row 23: 'PancakeSwap'
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.60'
$ws.Range('E23').Value = '  +4.21%  '

# Row 24: 'BitcoinCash'
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '265.68'
$ws.Range('E24').Value = '  -0.92%  '

# Row 25: 'ImmutableX'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.37'
$ws.Range('E25').Value = '  +7.21%  '

# Row 26: 'Dai'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '1.00'
$ws.Range('E26').Value = '  -0.68%  '

# Row 27: 'Filecoin'
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '7.72'
$ws.Range('E27').Value = '  +6.37%  '

# Row 28: 'Cosmos'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '11.41'
$ws.Range('E28').Value = '  +4.66%  '

# Row 29: 'Toncoin'
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.35'
$ws.Range('E29').Value = '  +2.43%  '

# Row 30: 'InjectiveProtocol'
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '40.26'
$ws.Range('E30').Value = '  +11.11%  '

# Row 31: 'Hedera'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0991'
$ws.Range('E31').Value = '  +15.41%  '

# Row 32: 'EthereumClassic'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '22.87'
$ws.Range('E32').Value = '  +2.46%  '

# Row 33: 'Monero'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '173.06'
$ws.Range('E33').Value = '  +4.95%  '

# Row 34: 'WEMIXToken'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '2.95'
$ws.Range('E34').Value = '  +12.06%  '

# Row 35: 'Stellar'
$ws.Range('B35').Value = 'Stellar'
$ws.Range('C35').Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.133'
$ws.Range('E35').Value = '  +1.99%  '

# Row 36: 'RenderToken'
$ws.Range('B36').Value = 'RenderToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.04'
$ws.Range('E36').Value = '  +11.02%  '

# Row 37: 'Kaspa'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.119'
$ws.Range('E37').Value = '  +6.76%  '

# Row 38: 'NEARProtocol'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '4.21'
$ws.Range('E38').Value = '  +16.06%  '

# Row 39: 'LidoDAOToken'
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '3.11'
$ws.Range('E39').Value = '  +10.28%  '

# Row 40: 'VeChain'
$ws.Range('E40').Value = '  +4.92%  '

# Row 41: 'ARBITRUM'
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.78'
$ws.Range('E41').Value = '  +11.89%  '

# Row 42: 'Celestia'
$ws.Range('B42').Value = 'Celestia'
$ws.Range('C42').Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '13.87'
$ws.Range('E42').Value = '  +12.47%  '

# Row 43: 'Algorand'
$ws.Range('B43').Value = 'Algorand'
$ws.Range('C43').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.243'
$ws.Range('E43').Value = '  +7.34%  '

# Row 44: 'BitcoinSV'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '100.34'
$ws.Range('E44').Value = '  -8.46%  '

# Row 45: 'MultiversX'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '72.33'
$ws.Range('E45').Value = '  +2.00%  '

# Row 46: 'ordi'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '89.75'
$ws.Range('E46').Value = '  +15.40%  '

# Row 47: 'FirstDigitalUSD'
$ws.Range('E47').Value = '  -0.50%  '

# Row 48: 'THORChain'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '5.87'
$ws.Range('E48').Value = '  +13.88%  '

# Row 49: 'Aave'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '116.51'
$ws.Range('E49').Value = '  +5.21%  '

# Row 50: 'FraxShare'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.48'
$ws.Range('E50').Value = '  +9.69%  '

# Row 51: 'MinaProtocolToken'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '1.58'
$ws.Range('E51').Value = '  +10.52%  '
